# Weekly update: insert the newest week's records (2 new rows) at the top of the
# Coliflor (Vega Monumental Concepción) data block, pushing the older rows down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 427-428; this shifts the existing rows 427-491 down to
# 429-493 (preserving all of their original values/styles), exactly matching the
# target edit.
$ws.Range("A427:A428").EntireRow.Insert()

# Row 427: new "Primera" quality record for the new reporting date (2023-10-19 ->
# Excel serial 45218).
$ws.Range("A427").Value2 = 11
$ws.Range("B427").Value2 = "Vega Monumental Concepción"
$ws.Range("C427").Value2 = "Bíobío"
$ws.Range("D427").Value2 = 45218
$ws.Range("E427").Value2 = 8
$ws.Range("F427").Value2 = 100112008
$ws.Range("G427").Value2 = "Coliflor"
$ws.Range("H427").Value2 = "Sin especificar"
$ws.Range("I427").Value2 = "Primera"
$ws.Range("J427").Value2 = 1000
$ws.Range("K427").Value2 = 1000
$ws.Range("L427").Value2 = 1000
$ws.Range("M427").Value2 = 1000
$ws.Range("N427").Value2 = "$/unidad"
$ws.Range("O427").Value2 = "Región Metropolitana"
$ws.Range("P427").Value2 = 1000
$ws.Range("Q427").Value2 = 1
$ws.Range("R427").Value2 = "Hortaliza"

# Row 428: new "Segunda" quality record for the same new reporting date.
$ws.Range("A428").Value2 = 11
$ws.Range("B428").Value2 = "Vega Monumental Concepción"
$ws.Range("C428").Value2 = "Bíobío"
$ws.Range("D428").Value2 = 45218
$ws.Range("E428").Value2 = 8
$ws.Range("F428").Value2 = 100112008
$ws.Range("G428").Value2 = "Coliflor"
$ws.Range("H428").Value2 = "Sin especificar"
$ws.Range("I428").Value2 = "Segunda"
$ws.Range("J428").Value2 = 1000
$ws.Range("K428").Value2 = 700
$ws.Range("L428").Value2 = 700
$ws.Range("M428").Value2 = 700
$ws.Range("N428").Value2 = "$/unidad"
$ws.Range("O428").Value2 = "Región Metropolitana"
$ws.Range("P428").Value2 = 700
$ws.Range("Q428").Value2 = 1
$ws.Range("R428").Value2 = "Hortaliza"
